# Applies the symbol-list update commit (2022-12-24 GitHub Actions run).
# - Column D ("Price") is refreshed with newly scraped quotes for most rows.
# - Column G ("Hora") moves from "4" to "5" for every data row (2-51).
# - The "Worstin24h" marker in column E moves from row 18 (One/ONE) to
#   row 48 (BOLO/BOLO), reflecting the new worst-24h-performer coin.
#
# A leading "'" is used on numeric-looking values so Excel stores them as
# literal text (matching the source file's inlineStr cells) instead of
# coercing them into real numbers and losing formatting such as
# trailing/leading zeros (e.g. "0.006200", "0.00000000751").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    ,@("D2", "'245.75")
    ,@("G2", "'5")
    ,@("D3", "'21.96")
    ,@("G3", "'5")
    ,@("D4", "'5.335")
    ,@("G4", "'5")
    ,@("D5", "'0.05954")
    ,@("G5", "'5")
    ,@("D6", "'3.397")
    ,@("G6", "'5")
    ,@("D7", "'6.393")
    ,@("G7", "'5")
    ,@("D8", "'0.8136")
    ,@("G8", "'5")
    ,@("D9", "'0.9642")
    ,@("G9", "'5")
    ,@("G10", "'5")
    ,@("D11", "'0.03703")
    ,@("G11", "'5")
    ,@("D12", "'0.07392")
    ,@("G12", "'5")
    ,@("D13", "'0.03046")
    ,@("G13", "'5")
    ,@("D14", "'0.09396")
    ,@("G14", "'5")
    ,@("D15", "'4.002")
    ,@("G15", "'5")
    ,@("D16", "'0.001597")
    ,@("G16", "'5")
    ,@("G17", "'5")
    ,@("E18", "17OneONE")
    ,@("G18", "'5")
    ,@("D19", "'0.006200")
    ,@("G19", "'5")
    ,@("D20", "'0.004143")
    ,@("G20", "'5")
    ,@("D21", "'0.0009860")
    ,@("G21", "'5")
    ,@("D22", "'0.00009707")
    ,@("G22", "'5")
    ,@("D23", "'3.742")
    ,@("G23", "'5")
    ,@("G24", "'5")
    ,@("G25", "'5")
    ,@("G26", "'5")
    ,@("G27", "'5")
    ,@("G28", "'5")
    ,@("G29", "'5")
    ,@("G30", "'5")
    ,@("G31", "'5")
    ,@("G32", "'5")
    ,@("G33", "'5")
    ,@("G34", "'5")
    ,@("G35", "'5")
    ,@("G36", "'5")
    ,@("G37", "'5")
    ,@("G38", "'5")
    ,@("G39", "'5")
    ,@("D40", "'0.03929")
    ,@("G40", "'5")
    ,@("D41", "'0.006525")
    ,@("G41", "'5")
    ,@("D42", "'0.1074")
    ,@("G42", "'5")
    ,@("D43", "'0.002702")
    ,@("G43", "'5")
    ,@("D44", "'0.005374")
    ,@("G44", "'5")
    ,@("D45", "'0.00005312")
    ,@("G45", "'5")
    ,@("D46", "'0.00000000751")
    ,@("G46", "'5")
    ,@("D47", "'0.8506")
    ,@("G47", "'5")
    ,@("D48", "'0.04127")
    ,@("E48", "47BOLOBOLOWorstin24h")
    ,@("G48", "'5")
    ,@("G49", "'5")
    ,@("G50", "'5")
    ,@("G51", "'5")
)

foreach ($update in $updates) {
    $ws.Range($update[0]).Value = $update[1]
}
